# Scheduled-runner style refresh of the market-board price/profit columns
# (H:N) on each class sheet -- values only, no structural changes.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 314.4
$ws.Range("I2").Value = 342.5
$ws.Range("K2").Value = 342.5
$ws.Range("M2").Value = -229.5
# Row 9
$ws.Range("H9").Value = 204.25
$ws.Range("J9").Value = 504.5
$ws.Range("L9").Value = 504.5
$ws.Range("N9").Value = -842.5
# Row 40
$ws.Range("H40").Value = 4323.4116
$ws.Range("I40").Value = 2499
$ws.Range("J40").Value = 4437.4375
$ws.Range("K40").Value = 2499
$ws.Range("L40").Value = 4437.4375
$ws.Range("N40").Value = -4787.4375
$ws.Range("M40").Value = -2324
# Row 53
$ws.Range("H53").Value = 284.5
$ws.Range("I53").Value = 281.30768
$ws.Range("J53").Value = 298.33334
$ws.Range("K53").Value = 281.30768
$ws.Range("L53").Value = 298.33334
$ws.Range("M53").Value = 355.69232
$ws.Range("N53").Value = -1572.33334
# Row 86
$ws.Range("H86").Value = 1049.5
$ws.Range("I86").Value = 999
$ws.Range("J86").Value = 1100
$ws.Range("K86").Value = 999
$ws.Range("L86").Value = 1100
$ws.Range("M86").Value = 124
$ws.Range("N86").Value = -3346
# Row 89
$ws.Range("H89").Value = 1049.5
$ws.Range("I89").Value = 999
$ws.Range("J89").Value = 1100
$ws.Range("K89").Value = 4995
$ws.Range("L89").Value = 5500
$ws.Range("M89").Value = 621
$ws.Range("N89").Value = -16732
# Row 107
$ws.Range("H107").Value = 1789
$ws.Range("I107").Value = 873
$ws.Range("K107").Value = 873
$ws.Range("M107").Value = 1047

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 811
$ws.Range("I2").Value = 811
$ws.Range("K2").Value = 811
$ws.Range("M2").Value = -698
# Row 5
$ws.Range("H5").Value = 16.777779
$ws.Range("J5").Value = 18.571428
$ws.Range("L5").Value = 18.571428
$ws.Range("N5").Value = -242.571428
# Row 32
$ws.Range("H32").Value = 3890.9092
$ws.Range("I32").Value = 4180
$ws.Range("K32").Value = 4180
$ws.Range("M32").Value = -3893
# Row 97
$ws.Range("H97").Value = 2283.5
$ws.Range("I97").Value = 1407.6666
$ws.Range("K97").Value = 1407.6666
$ws.Range("M97").Value = -911.6666
# Row 116
$ws.Range("H116").Value = 811
$ws.Range("I116").Value = 811
$ws.Range("K116").Value = 811
$ws.Range("M116").Value = 1483

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 811
$ws.Range("I3").Value = 811
$ws.Range("K3").Value = 811
$ws.Range("M3").Value = -697
# Row 4
$ws.Range("H4").Value = 16.777779
$ws.Range("J4").Value = 18.571428
$ws.Range("L4").Value = 18.571428
$ws.Range("N4").Value = -248.571428
# Row 22
$ws.Range("H22").Value = 952
$ws.Range("I22").Value = 1014.2308
$ws.Range("J22").Value = 749.75
$ws.Range("K22").Value = 1014.2308
$ws.Range("L22").Value = 749.75
$ws.Range("M22").Value = -841.2308
$ws.Range("N22").Value = -1095.75
# Row 99
$ws.Range("H99").Value = 200
$ws.Range("I99").Value = 200
$ws.Range("K99").Value = 200
$ws.Range("M99").Value = 1298
# Row 105
$ws.Range("H105").Value = 5966.6665
$ws.Range("I105").Value = 8000
$ws.Range("J105").Value = 4950
$ws.Range("K105").Value = 8000
$ws.Range("L105").Value = 4950
$ws.Range("M105").Value = -6253
$ws.Range("N105").Value = -8444
# Row 116
$ws.Range("H116").Value = 92497.5
$ws.Range("J116").Value = 92497.5
$ws.Range("L116").Value = 92497.5
$ws.Range("N116").Value = -101675.5

$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
# Row 7
$ws.Range("H7").Value = 87.916664
$ws.Range("I7").Value = 87.25
$ws.Range("J7").Value = 89.25
$ws.Range("K7").Value = 87.25
$ws.Range("L7").Value = 89.25
$ws.Range("M7").Value = 25.75
$ws.Range("N7").Value = -315.25
# Row 31
$ws.Range("H31").Value = 4453.5713
$ws.Range("I31").Value = 1870.5
$ws.Range("K31").Value = 1870.5
$ws.Range("M31").Value = -1575.5
# Row 34
$ws.Range("H34").Value = 4453.5713
$ws.Range("I34").Value = 1870.5
$ws.Range("K34").Value = 1870.5
$ws.Range("M34").Value = -1668.5
# Row 58
$ws.Range("H58").Value = 5000
$ws.Range("I58").Value = 5000
$ws.Range("K58").Value = 5000
$ws.Range("M58").Value = -4797
# Row 94
$ws.Range("H94").Value = 1775
$ws.Range("I94").Value = 100
$ws.Range("J94").Value = 2333.3333
$ws.Range("K94").Value = 100
$ws.Range("L94").Value = 2333.3333
$ws.Range("N94").Value = -3235.3333
$ws.Range("M94").Value = 351
# Row 134
$ws.Range("H134").Value = 1096.25
$ws.Range("I134").Value = 1099
$ws.Range("J134").Value = 1093.5
$ws.Range("K134").Value = 3297
$ws.Range("L134").Value = 3280.5
$ws.Range("M134").Value = -762
$ws.Range("N134").Value = -8350.5
# Row 136
$ws.Range("H136").Value = 5000
$ws.Range("I136").Value = 5000
$ws.Range("K136").Value = 15000
$ws.Range("M136").Value = -12450

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 462007.47
$ws.Range("I4").Value = 417174.84
$ws.Range("J4").Value = 999999
$ws.Range("K4").Value = 1251524.52
$ws.Range("L4").Value = 2999997
$ws.Range("M4").Value = -1251412.52
$ws.Range("N4").Value = -3000221
# Row 129
$ws.Range("H129").Value = 1559
$ws.Range("I129").Value = 2133.3333
$ws.Range("J129").Value = 697.5
$ws.Range("K129").Value = 6399.999899999999
$ws.Range("L129").Value = 2092.5
$ws.Range("M129").Value = -1399.999899999999
$ws.Range("N129").Value = -12092.5
# Row 131
$ws.Range("H131").Value = 1928.625
$ws.Range("J131").Value = 3458.3333
$ws.Range("L131").Value = 10374.9999
$ws.Range("N131").Value = -20454.9999

$ws = $wb.Worksheets.Item("GSM")
# Row 11
$ws.Range("H11").Value = 4821832
$ws.Range("J11").Value = 551
$ws.Range("L11").Value = 551
$ws.Range("N11").Value = -829
# Row 18
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()
# Row 53
$ws.Range("H53").Value = 15000
$ws.Range("J53").Value = 15000
$ws.Range("L53").Value = 15000
$ws.Range("N53").Value = -16262
# Row 97
$ws.Range("H97").Value = 1925
$ws.Range("I97").Value = 1925
$ws.Range("K97").Value = 1925
$ws.Range("M97").Value = -1429

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 2067.9375
$ws.Range("I22").Value = 1826.091
$ws.Range("K22").Value = 1826.091
$ws.Range("M22").Value = -1531.091
# Row 24
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("M24").ClearContents()
$ws.Range("N24").ClearContents()
# Row 27
$ws.Range("H27").Value = 2067.9375
$ws.Range("I27").Value = 1826.091
$ws.Range("K27").Value = 1826.091
$ws.Range("M27").Value = -1719.091
# Row 43
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
# Row 46
$ws.Range("H46").Value = 4668.683
$ws.Range("J46").Value = 4822.054
$ws.Range("L46").Value = 4822.054
$ws.Range("N46").Value = -5198.054
# Row 68
$ws.Range("H68").Value = 3149
$ws.Range("I68").Value = 2753.4546
$ws.Range("K68").Value = 2753.4546
$ws.Range("M68").Value = -2004.4546
# Row 71
$ws.Range("H71").Value = 3149
$ws.Range("I71").Value = 2753.4546
$ws.Range("K71").Value = 13767.273
$ws.Range("M71").Value = -10023.273
# Row 132
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()
# Row 136
$ws.Range("H136").Value = 6666.3335
$ws.Range("I136").Value = 8000
$ws.Range("K136").Value = 24000
$ws.Range("M136").Value = -21450

$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 6800000
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
# Row 20
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
# Row 25
$ws.Range("H25").Value = 28397.5
$ws.Range("J25").Value = 28397.5
$ws.Range("L25").Value = 28397.5
$ws.Range("N25").Value = -28983.5
# Row 32
$ws.Range("H32").Value = 12013
$ws.Range("I32").Value = 12013
$ws.Range("K32").Value = 12013
$ws.Range("M32").Value = -11696
# Row 132
$ws.Range("H132").Value = 1000
$ws.Range("I132").Value = 1000
$ws.Range("K132").Value = 3000
$ws.Range("M132").Value = -470
# Row 136
$ws.Range("H136").Value = 4544.6665
$ws.Range("J136").Value = 10000
$ws.Range("L136").Value = 30000
$ws.Range("N136").Value = -35100
